$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume/1h change (E) columns
# D-column values are forced to Text format while assigning so that
# numeric-looking strings (e.g. "306.15", "35.00") keep their exact
# textual representation instead of being auto-converted to numbers,
# then the style is reset back to Normal so no stray formatting is left
# behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.933.52'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.334.52'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.510'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.61%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.20'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0798'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.86'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.368.08'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.811'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.868.66'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0910'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.69'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.48%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.94'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.45%  '
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.03'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.38%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.88%  '
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.48'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.02%  '
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.92'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.62%  '
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('E40').Value = '  -2.68%  '
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.48'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.004.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.65%  '
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.02%  '
$ws.Range('E47').Value = '  -1.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.62'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.561.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.24%  '
